$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  ,@(8,'Terminal La Palmera de La Serena','Coquimbo',44858,4,100112013,'Alcachofa','Española','Primera',500,9500,10000,9750,'$/caja 30 unidades','Provincia del Elquí',325,30,'Hortaliza')
  ,@(8,'Terminal La Palmera de La Serena','Coquimbo',45154,4,100112013,'Alcachofa','Española','Primera',500,11000,12000,11500,'$/caja 30 unidades','Provincia de Limarí',383,30,'Hortaliza')
  ,@(8,'Terminal La Palmera de La Serena','Coquimbo',45037,4,100112013,'Alcachofa','Madrigal','Primera',400,16000,17000,16500,'$/caja 40 unidades','Provincia del Elquí',412,40,'Hortaliza')
  ,@(8,'Terminal La Palmera de La Serena','Coquimbo',45124,4,100112013,'Alcachofa','Española','Primera',400,14000,15000,14500,'$/caja 30 unidades','Provincia de Limarí',483,30,'Hortaliza')
  ,@(8,'Terminal La Palmera de La Serena','Coquimbo',45069,4,100112013,'Alcachofa','Madrigal','Primera',500,17000,18000,17500,'$/caja 40 unidades','Provincia del Elquí',438,40,'Hortaliza')
  ,@(8,'Terminal La Palmera de La Serena','Coquimbo',45126,4,100112013,'Alcachofa','Española','Primera',500,14000,15000,14500,'$/caja 30 unidades','Provincia de Limarí',483,30,'Hortaliza')
  ,@(8,'Terminal La Palmera de La Serena','Coquimbo',45135,4,100112013,'Alcachofa','Madrigal','Primera',500,10000,11000,10500,'$/caja 40 unidades','Provincia del Elquí',262,40,'Hortaliza')
  ,@(8,'Terminal La Palmera de La Serena','Coquimbo',44855,4,100112013,'Alcachofa','Española','Primera',540,9500,10000,9750,'$/caja 30 unidades','Provincia del Elquí',325,30,'Hortaliza')
  ,@(8,'Terminal La Palmera de La Serena','Coquimbo',45118,4,100112013,'Alcachofa','Española','Primera',360,14000,15000,14500,'$/caja 30 unidades','Provincia de Limarí',483,30,'Hortaliza')
  ,@(8,'Terminal La Palmera de La Serena','Coquimbo',45118,4,100112013,'Alcachofa','Madrigal','Primera',400,12500,13000,12750,'$/caja 40 unidades','Provincia del Elquí',319,40,'Hortaliza')
  ,@(8,'Terminal La Palmera de La Serena','Coquimbo',44701,4,100112013,'Alcachofa','Española','Primera',400,19000,20000,19500,'$/caja 30 unidades','Provincia del Elquí',650,30,'Hortaliza')
  ,@(8,'Terminal La Palmera de La Serena','Coquimbo',44427,4,100112013,'Alcachofa','Madrigal','Primera',400,12000,13000,12500,'$/caja 40 unidades','Provincia de Limarí',312,40,'Hortaliza')
  ,@(8,'Terminal La Palmera de La Serena','Coquimbo',44839,4,100112013,'Alcachofa','Española','Primera',400,12000,13000,12500,'$/caja 30 unidades','Provincia del Elquí',417,30,'Hortaliza')
  ,@(8,'Terminal La Palmera de La Serena','Coquimbo',45204,4,100112013,'Alcachofa','Española','Primera',500,8000,9000,8500,'$/caja 30 unidades','Provincia de Limarí',283,30,'Hortaliza')
  ,@(8,'Terminal La Palmera de La Serena','Coquimbo',45082,4,100112013,'Alcachofa','Madrigal','Primera',200,17000,18000,17500,'$/caja 40 unidades','Provincia del Elquí',438,40,'Hortaliza')
  ,@(8,'Terminal La Palmera de La Serena','Coquimbo',45070,4,100112013,'Alcachofa','Madrigal','Primera',360,17000,18000,17500,'$/caja 40 unidades','Provincia del Elquí',438,40,'Hortaliza')
  ,@(8,'Terminal La Palmera de La Serena','Coquimbo',45090,4,100112013,'Alcachofa','Madrigal','Primera',340,15500,16000,15750,'$/caja 40 unidades','Provincia del Elquí',394,40,'Hortaliza')
  ,@(8,'Terminal La Palmera de La Serena','Coquimbo',44767,4,100112013,'Alcachofa','Madrigal','Primera',600,14000,15000,14500,'$/caja 40 unidades','Provincia de Limarí',362,40,'Hortaliza')
  ,@(8,'Terminal La Palmera de La Serena','Coquimbo',44426,4,100112013,'Alcachofa','Española','Primera',600,11500,12000,11750,'$/caja 30 unidades','Provincia de Limarí',392,30,'Hortaliza')
  ,@(8,'Terminal La Palmera de La Serena','Coquimbo',44426,4,100112013,'Alcachofa','Madrigal','Primera',500,12500,13000,12750,'$/caja 40 unidades','Provincia de Limarí',319,40,'Hortaliza')
  ,@(8,'Terminal La Palmera de La Serena','Coquimbo',44484,4,100112013,'Alcachofa','Española','Primera',300,9000,10000,9500,'$/caja 30 unidades','Provincia del Elquí',317,30,'Hortaliza')
  ,@(8,'Terminal La Palmera de La Serena','Coquimbo',44498,4,100112013,'Alcachofa','Española','Primera',400,8500,9000,8750,'$/caja 30 unidades','Provincia de Limarí',292,30,'Hortaliza')
  ,@(8,'Terminal La Palmera de La Serena','Coquimbo',45128,4,100112013,'Alcachofa','Española','Primera',400,14000,15000,14500,'$/caja 30 unidades','Provincia de Limarí',483,30,'Hortaliza')
  ,@(8,'Terminal La Palmera de La Serena','Coquimbo',45223,4,100112013,'Alcachofa','Española','Primera',440,9000,10000,9500,'$/caja 30 unidades','Provincia de Limarí',317,30,'Hortaliza')
  ,@(8,'Terminal La Palmera de La Serena','Coquimbo',45163,4,100112013,'Alcachofa','Española','Primera',600,10000,10500,10250,'$/caja 30 unidades','Provincia del Elquí',342,30,'Hortaliza')
  ,@(8,'Terminal La Palmera de La Serena','Coquimbo',44687,4,100112013,'Alcachofa','Española','Primera',400,18000,19000,18500,'$/caja 30 unidades','Provincia de Limarí',617,30,'Hortaliza')
  ,@(8,'Terminal La Palmera de La Serena','Coquimbo',45205,4,100112013,'Alcachofa','Española','Primera',520,8000,9000,8500,'$/caja 30 unidades','Provincia de Limarí',283,30,'Hortaliza')
  ,@(8,'Terminal La Palmera de La Serena','Coquimbo',45079,4,100112013,'Alcachofa','Madrigal','Primera',400,16500,17000,16750,'$/caja 40 unidades','Provincia del Elquí',419,40,'Hortaliza')
  ,@(8,'Terminal La Palmera de La Serena','Coquimbo',45180,4,100112013,'Alcachofa','Española','Primera',400,9000,10000,9500,'$/caja 30 unidades','Provincia del Elquí',317,30,'Hortaliza')
  ,@(8,'Terminal La Palmera de La Serena','Coquimbo',44420,4,100112013,'Alcachofa','Madrigal','Primera',800,14000,15000,14500,'$/caja 40 unidades','Provincia de Limarí',362,40,'Hortaliza')
  ,@(8,'Terminal La Palmera de La Serena','Coquimbo',44420,4,100112013,'Alcachofa','Madrigal','Primera',700,13000,14000,13500,'$/caja 40 unidades','Provincia del Elquí',338,40,'Hortaliza')
  ,@(8,'Terminal La Palmera de La Serena','Coquimbo',45142,4,100112013,'Alcachofa','Española','Primera',560,14000,15000,14500,'$/caja 30 unidades','Provincia del Elquí',483,30,'Hortaliza')
  ,@(8,'Terminal La Palmera de La Serena','Coquimbo',45142,4,100112013,'Alcachofa','Madrigal','Primera',600,9000,10000,9500,'$/caja 40 unidades','Provincia del Elquí',238,40,'Hortaliza')
  ,@(8,'Terminal La Palmera de La Serena','Coquimbo',45194,4,100112013,'Alcachofa','Española','Primera',600,8500,9000,8750,'$/caja 30 unidades','Provincia de Limarí',292,30,'Hortaliza')
  ,@(8,'Terminal La Palmera de La Serena','Coquimbo',44784,4,100112013,'Alcachofa','Madrigal','Primera',520,11500,12000,11750,'$/caja 40 unidades','Provincia del Elquí',294,40,'Hortaliza')
  ,@(8,'Terminal La Palmera de La Serena','Coquimbo',44438,4,100112013,'Alcachofa','Española','Primera',400,11000,12000,11500,'$/caja 30 unidades','Provincia del Elquí',383,30,'Hortaliza')
  ,@(8,'Terminal La Palmera de La Serena','Coquimbo',44729,4,100112013,'Alcachofa','Madrigal','Primera',400,16000,17000,16500,'$/caja 40 unidades','Provincia del Elquí',412,40,'Hortaliza')
  ,@(8,'Terminal La Palmera de La Serena','Coquimbo',44790,4,100112013,'Alcachofa','Española','Primera',400,14000,15000,14500,'$/caja 30 unidades','Provincia de Limarí',483,30,'Hortaliza')
  ,@(8,'Terminal La Palmera de La Serena','Coquimbo',44790,4,100112013,'Alcachofa','Madrigal','Primera',400,11500,12000,11750,'$/caja 40 unidades','Provincia del Elquí',294,40,'Hortaliza')
  ,@(8,'Terminal La Palmera de La Serena','Coquimbo',45224,4,100112013,'Alcachofa','Española','Primera',520,9000,10000,9500,'$/caja 30 unidades','Provincia de Limarí',317,30,'Hortaliza')
  ,@(8,'Terminal La Palmera de La Serena','Coquimbo',45134,4,100112013,'Alcachofa','Española','Primera',500,12500,13000,12750,'$/caja 30 unidades','Provincia de Limarí',425,30,'Hortaliza')
  ,@(8,'Terminal La Palmera de La Serena','Coquimbo',45132,4,100112013,'Alcachofa','Española','Primera',420,13000,14000,13500,'$/caja 30 unidades','Provincia de Limarí',450,30,'Hortaliza')
  ,@(8,'Terminal La Palmera de La Serena','Coquimbo',45114,4,100112013,'Alcachofa','Española','Primera',400,14000,15000,14500,'$/caja 30 unidades','Provincia de Limarí',483,30,'Hortaliza')
  ,@(8,'Terminal La Palmera de La Serena','Coquimbo',45114,4,100112013,'Alcachofa','Madrigal','Primera',400,15000,16000,15500,'$/caja 40 unidades','Provincia del Elquí',388,40,'Hortaliza')
  ,@(8,'Terminal La Palmera de La Serena','Coquimbo',45177,4,100112013,'Alcachofa','Española','Primera',600,7500,8000,7750,'$/caja 30 unidades','Provincia de Limarí',258,30,'Hortaliza')
  ,@(8,'Terminal La Palmera de La Serena','Coquimbo',45049,4,100112013,'Alcachofa','Madrigal','Primera',500,17000,18000,17500,'$/caja 40 unidades','Provincia del Elquí',438,40,'Hortaliza')
)

for ($i = 0; $i -lt $data.Count; $i++) {
  $rowNum = $i + 2
  $rowVals = $data[$i]
  for ($j = 0; $j -lt $rowVals.Count; $j++) {
    $colNum = $j + 1
    $ws.Cells.Item($rowNum, $colNum).Value = $rowVals[$j]
  }
}

"Applied permutation of rows 2-47 (weekly roll)."